$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Windows")

# Rows 10-13 describe the PXE (OS automated install) test programs.
# Commit: the PXE program's files/SOP+xml have now been committed, so:
#   - "Program Folder" (column C) moves from the placeholder "n/a" to the
#     real folder name "pxe"
#   - "Development Status" (column G) moves from "Processing" to "Delivered"
#     (copy G2's cell format, which already carries the green "Delivered"
#     style, onto G10:G13 so the color/border match the rest of the sheet)
$ws.Range("C10").Value = "pxe"
$ws.Range("C11").Value = "pxe"
$ws.Range("C12").Value = "pxe"
$ws.Range("C13").Value = "pxe"

$ws.Range("G2").Copy() | Out-Null
$ws.Range("G10:G13").PasteSpecial(-4122) | Out-Null

$ws.Range("G10").Value = "Delivered"
$ws.Range("G11").Value = "Delivered"
$ws.Range("G12").Value = "Delivered"
$ws.Range("G13").Value = "Delivered"

$excel.CutCopyMode = $false
